$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-24 20:56:11"
$wsZhCn.Range("H4").Value = "2016-03-24 20:56:52"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-24 20:56:17"
$wsDeDe.Range("H4").Value = "2016-03-24 20:56:59"
